# Update the Ylo column (D) for all stimuli rows (2-121) from 1 to 200.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D121").Value = 200

# Move the active selection to match the saved view state (F7).
[void]$ws.Range("F7").Select()
